$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1100634.9
$ws.Range("I41").Value = 1429193.1
$ws.Range("J41").Value = 333999
$ws.Range("K41").Value = 1429193.1
$ws.Range("L41").Value = 333999
$ws.Range("M41").Value = -1428753.1
$ws.Range("N41").Value = -334879

$ws.Range("H46").Value = 9333.546
$ws.Range("I46").Value = 787.5
$ws.Range("J46").Value = 14217
$ws.Range("K46").Value = 2362.5
$ws.Range("L46").Value = 42651
$ws.Range("M46").Value = -2243.5
$ws.Range("N46").Value = -42889

$ws.Range("H60").Value = 9333.546
$ws.Range("I60").Value = 787.5
$ws.Range("J60").Value = 14217
$ws.Range("K60").Value = 2362.5
$ws.Range("L60").Value = 42651
$ws.Range("M60").Value = -1878.5
$ws.Range("N60").Value = -43619

$ws.Range("H62").Value = 704566.9
$ws.Range("I62").Value = 1669666.4
$ws.Range("J62").Value = 61167.223
$ws.Range("K62").Value = 1669666.4
$ws.Range("L62").Value = 61167.223
$ws.Range("M62").Value = -1669042.4
$ws.Range("N62").Value = -62415.223

$ws.Range("H64").Value = 3930.0667
$ws.Range("I64").Value = 3661.2222
$ws.Range("J64").Value = 4333.3335
$ws.Range("K64").Value = 3661.2222
$ws.Range("L64").Value = 4333.3335
$ws.Range("M64").Value = -3413.2222
$ws.Range("N64").Value = -4829.3335

$ws.Range("H65").Value = 704566.9
$ws.Range("I65").Value = 1669666.4
$ws.Range("J65").Value = 61167.223
$ws.Range("K65").Value = 8348332
$ws.Range("L65").Value = 305836.115
$ws.Range("M65").Value = -8345212
$ws.Range("N65").Value = -312076.115

$ws.Range("H67").Value = 3930.0667
$ws.Range("I67").Value = 3661.2222
$ws.Range("J67").Value = 4333.3335
$ws.Range("K67").Value = 3661.2222
$ws.Range("L67").Value = 4333.3335
$ws.Range("M67").Value = -2803.2222
$ws.Range("N67").Value = -6049.3335

$ws.Range("H86").Value = 14362336
$ws.Range("I86").Value = 3501.5
$ws.Range("J86").Value = 16755475
$ws.Range("K86").Value = 3501.5
$ws.Range("L86").Value = 16755475
$ws.Range("M86").Value = -2378.5
$ws.Range("N86").Value = -16757721

$ws.Range("H88").Value = 2447.2222
$ws.Range("I88").Value = 2500.75
$ws.Range("J88").Value = 2431.9285
$ws.Range("K88").Value = 2500.75
$ws.Range("L88").Value = 2431.9285
$ws.Range("M88").Value = -2094.75
$ws.Range("N88").Value = -3243.9285

$ws.Range("H89").Value = 14362336
$ws.Range("I89").Value = 3501.5
$ws.Range("J89").Value = 16755475
$ws.Range("K89").Value = 17507.5
$ws.Range("L89").Value = 83777375
$ws.Range("M89").Value = -11891.5
$ws.Range("N89").Value = -83788607

$ws.Range("H91").Value = 2447.2222
$ws.Range("I91").Value = 2500.75
$ws.Range("J91").Value = 2431.9285
$ws.Range("K91").Value = 2500.75
$ws.Range("L91").Value = 2431.9285
$ws.Range("M91").Value = -1096.75
$ws.Range("N91").Value = -5239.9285

$ws.Range("H98").Value = 2537.7036
$ws.Range("I98").Value = 2648.2727
$ws.Range("J98").Value = 2051.2
$ws.Range("K98").Value = 2648.2727
$ws.Range("L98").Value = 2051.2
$ws.Range("M98").Value = -1150.2727
$ws.Range("N98").Value = -5047.2

$ws.Range("H106").Value = 72623.5
$ws.Range("I106").Value = 1227.4166
$ws.Range("J106").Value = 501000
$ws.Range("K106").Value = 1227.4166
$ws.Range("L106").Value = 501000
$ws.Range("M106").Value = -596.4166
$ws.Range("N106").Value = -502262

$ws.Range("H122").Value = 2537.7036
$ws.Range("I122").Value = 2648.2727
$ws.Range("J122").Value = 2051.2
$ws.Range("K122").Value = 7944.8181
$ws.Range("L122").Value = 6153.599999999999
$ws.Range("M122").Value = -5494.8181
$ws.Range("N122").Value = -11053.6

$ws.Range("H132").Value = 3055.625
$ws.Range("I132").Value = 2847.0815
$ws.Range("J132").Value = 4515.4287
$ws.Range("K132").Value = 8541.244499999999
$ws.Range("L132").Value = 13546.2861
$ws.Range("M132").Value = -6011.244499999999
$ws.Range("N132").Value = -18606.2861

$ws.Range("H137").Value = 2328.1304
$ws.Range("I137").Value = 2267.85
$ws.Range("J137").Value = 2730
$ws.Range("K137").Value = 6803.549999999999
$ws.Range("L137").Value = 8190
$ws.Range("M137").Value = -4253.549999999999
$ws.Range("N137").Value = -13290

$ws.Range("H138").Value = 1788.1642
$ws.Range("I138").Value = 1031.317
$ws.Range("J138").Value = 2981.6538
$ws.Range("K138").Value = 3093.951
$ws.Range("L138").Value = 8944.9614
$ws.Range("M138").Value = 2046.049
$ws.Range("N138").Value = -19224.9614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2829.7856
$ws.Range("I45").Value = 2570.4
$ws.Range("J45").Value = 3478.25
$ws.Range("K45").Value = 2570.4
$ws.Range("L45").Value = 3478.25
$ws.Range("M45").Value = -2193.4
$ws.Range("N45").Value = -4232.25

$ws.Range("H74").Value = 5153.095
$ws.Range("I74").Value = 5874.75
$ws.Range("J74").Value = 2843.8
$ws.Range("K74").Value = 5874.75
$ws.Range("L74").Value = 2843.8
$ws.Range("M74").Value = -5000.75
$ws.Range("N74").Value = -4591.8

$ws.Range("H77").Value = 5153.095
$ws.Range("I77").Value = 5874.75
$ws.Range("J77").Value = 2843.8
$ws.Range("K77").Value = 29373.75
$ws.Range("L77").Value = 14219
$ws.Range("M77").Value = -25005.75
$ws.Range("N77").Value = -22955

$ws.Range("H132").Value = 3952.8
$ws.Range("I132").Value = 3566.125
$ws.Range("J132").Value = 5499.5
$ws.Range("K132").Value = 10698.375
$ws.Range("L132").Value = 16498.5
$ws.Range("M132").Value = -8168.375
$ws.Range("N132").Value = -21558.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1467
$ws.Range("I134").Value = 1467
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4401
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1866

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 4598
$ws.Range("I3").Value = 5500
$ws.Range("J3").Value = 990
$ws.Range("K3").Value = 5500
$ws.Range("L3").Value = 990
$ws.Range("M3").Value = -5387
$ws.Range("N3").Value = -1216

$ws.Range("H22").Value = 278.33334
$ws.Range("I22").Value = 278.33334
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 278.33334
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 71.66665999999998

$ws.Range("H31").Value = 4570.38
$ws.Range("I31").Value = 15116.556
$ws.Range("J31").Value = 3214.4429
$ws.Range("K31").Value = 15116.556
$ws.Range("L31").Value = 3214.4429
$ws.Range("M31").Value = -14821.556
$ws.Range("N31").Value = -3804.4429

$ws.Range("H34").Value = 4570.38
$ws.Range("I34").Value = 15116.556
$ws.Range("J34").Value = 3214.4429
$ws.Range("K34").Value = 15116.556
$ws.Range("L34").Value = 3214.4429
$ws.Range("M34").Value = -14914.556
$ws.Range("N34").Value = -3618.4429

$ws.Range("H122").Value = 2351.6667
$ws.Range("I122").Value = 2672
$ws.Range("J122").Value = 750
$ws.Range("K122").Value = 8016
$ws.Range("L122").Value = 2250
$ws.Range("M122").Value = -5566
$ws.Range("N122").Value = -7150

$ws.Range("H132").Value = 371451.16
$ws.Range("I132").Value = 228351.94
$ws.Range("J132").Value = 1001087.8
$ws.Range("K132").Value = 685055.8200000001
$ws.Range("L132").Value = 3003263.4
$ws.Range("M132").Value = -682525.8200000001
$ws.Range("N132").Value = -3008323.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = $null

$ws.Range("H87").Value = 95
$ws.Range("I87").Value = 95
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 285
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = 963

$ws.Range("H90").Value = 95
$ws.Range("I90").Value = 95
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 855
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = 5385

$ws.Range("H92").Value = 321.66666
$ws.Range("I92").Value = 356.66666
$ws.Range("J92").Value = 304.16666
$ws.Range("K92").Value = 1069.99998
$ws.Range("L92").Value = 912.4999799999999
$ws.Range("M92").Value = 178.0000199999999
$ws.Range("N92").Value = -3408.49998

$ws.Range("H107").Value = 623.4666999999999
$ws.Range("I107").Value = 618
$ws.Range("J107").Value = 700
$ws.Range("K107").Value = 1854
$ws.Range("L107").Value = 2100
$ws.Range("M107").Value = 66
$ws.Range("N107").Value = -5940

$ws.Range("H117").Value = 97.75
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 97.75
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 293.25
$ws.Range("N117").Value = -7177.25
$ws.Range("M117").Value = $null

$ws.Range("H132").Value = 2066.4243
$ws.Range("I132").Value = 1828.0952
$ws.Range("J132").Value = 2483.5
$ws.Range("K132").Value = 16452.8568
$ws.Range("L132").Value = 22351.5
$ws.Range("M132").Value = -13922.8568
$ws.Range("N132").Value = -27411.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 278.18182
$ws.Range("I2").Value = 24.333334
$ws.Range("J2").Value = 582.8
$ws.Range("K2").Value = 24.333334
$ws.Range("L2").Value = 582.8
$ws.Range("M2").Value = 88.66666599999999
$ws.Range("N2").Value = -808.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 668.1
$ws.Range("I22").Value = 420.25
$ws.Range("J22").Value = 833.3333
$ws.Range("K22").Value = 420.25
$ws.Range("L22").Value = 833.3333
$ws.Range("M22").Value = -125.25
$ws.Range("N22").Value = -1423.3333

$ws.Range("H27").Value = 668.1
$ws.Range("I27").Value = 420.25
$ws.Range("J27").Value = 833.3333
$ws.Range("K27").Value = 420.25
$ws.Range("L27").Value = 833.3333
$ws.Range("M27").Value = -313.25
$ws.Range("N27").Value = -1047.3333

$ws.Range("H40").Value = 2453
$ws.Range("I40").Value = 2238.8
$ws.Range("J40").Value = 2988.5
$ws.Range("K40").Value = 2238.8
$ws.Range("L40").Value = 2988.5
$ws.Range("M40").Value = -2102.8
$ws.Range("N40").Value = -3260.5

$ws.Range("H68").Value = 11741.875
$ws.Range("I68").Value = 3237
$ws.Range("J68").Value = 25916.666
$ws.Range("K68").Value = 3237
$ws.Range("L68").Value = 25916.666
$ws.Range("M68").Value = -2488
$ws.Range("N68").Value = -27414.666

$ws.Range("H71").Value = 11741.875
$ws.Range("I71").Value = 3237
$ws.Range("J71").Value = 25916.666
$ws.Range("K71").Value = 16185
$ws.Range("L71").Value = 129583.33
$ws.Range("M71").Value = -12441
$ws.Range("N71").Value = -137071.33

$ws.Range("H122").Value = 3781.9285
$ws.Range("I122").Value = 2703.9167
$ws.Range("J122").Value = 10250
$ws.Range("K122").Value = 8111.750100000001
$ws.Range("L122").Value = 30750
$ws.Range("M122").Value = -5661.750100000001
$ws.Range("N122").Value = -35650

$ws.Range("H132").Value = 43481870
$ws.Range("I132").Value = 52635012
$ws.Range("J132").Value = 4451
$ws.Range("K132").Value = 157905036
$ws.Range("L132").Value = 13353
$ws.Range("M132").Value = -157902506
$ws.Range("N132").Value = -18413

$ws.Range("H136").Value = 2063.9607
$ws.Range("I136").Value = 1699.7778
$ws.Range("J136").Value = 2938
$ws.Range("K136").Value = 5099.3334
$ws.Range("L136").Value = 8814
$ws.Range("M136").Value = -2549.3334
$ws.Range("N136").Value = -13914

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 7316
$ws.Range("I113").Value = 7269.6665
$ws.Range("J113").Value = 7489.75
$ws.Range("K113").Value = 21808.9995
$ws.Range("L113").Value = 22469.25
$ws.Range("M113").Value = -19638.9995
$ws.Range("N113").Value = -26809.25

$ws.Range("H132").Value = 381217.06
$ws.Range("I132").Value = 428980.97
$ws.Range("J132").Value = 7066.6665
$ws.Range("K132").Value = 1286942.91
$ws.Range("L132").Value = 24600
$ws.Range("M132").Value = -1284412.91
$ws.Range("N132").Value = -29660

$ws.Range("H136").Value = 6609.3076
$ws.Range("I136").Value = 3761.3928
$ws.Range("J136").Value = 13858.546
$ws.Range("K136").Value = 11284.1784
$ws.Range("L136").Value = 41575.638
$ws.Range("M136").Value = -8734.178400000001
$ws.Range("N136").Value = -46675.638
